# Applies the "Updated cryptos list" crypto-price/volume refresh.
# Source: cryptos.xlsx row-level diff (rows 2-51, columns B-E).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '95.164.12'
$ws.Range("E2").Value = '  -1.67%  '
$ws.Range("D3").Value = '3.492.04'
$ws.Range("E3").Value = '  +4.54%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.16'
$ws.Range("E5").Value = '  -4.17%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '646.66'
$ws.Range("E6").Value = '  -1.47%  '
$ws.Range("E7").Value = '  +4.21%  '
$ws.Range("E8").Value = '  -3.90%  '
$ws.Range("E9").Value = '  +0.06%  '
$ws.Range("E10").Value = '  +0.48%  '
$ws.Range("D11").Value = '3.489.36'
$ws.Range("E11").Value = '  +4.53%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '43.13'
$ws.Range("E12").Value = '  +6.15%  '
$ws.Range("E13").Value = '  -3.70%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.22'
$ws.Range("E14").Value = '  +2.02%  '
$ws.Range("D15").Value = '95.071.50'
$ws.Range("E15").Value = '  -1.45%  '
$ws.Range("D16").Value = '4.149.58'
$ws.Range("E16").Value = '  +4.62%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000258'
$ws.Range("E17").Value = '  +2.31%  '
$ws.Range("E18").Value = '  -2.99%  '
$ws.Range("D19").Value = '3.496.70'
$ws.Range("E19").Value = '  +4.33%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.11'
$ws.Range("E20").Value = '  +4.35%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.50'
$ws.Range("E21").Value = '  +7.94%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.514'
$ws.Range("E22").Value = '  -10.86%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '508.01'
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.22'
$ws.Range("E24").Value = '  -4.27%  '
$ws.Range("B25").Value = 'NEARProtocol'
$ws.Range("C25").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '6.72'
$ws.Range("E25").Value = '  +2.02%  '
$ws.Range("B26").Value = 'PEPE'
$ws.Range("C26").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000194'
$ws.Range("E26").Value = '  -2.39%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '95.63'
$ws.Range("E27").Value = '  -0.99%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.25'
$ws.Range("E28").Value = '  +0.95%  '
$ws.Range("B29").Value = 'InternetComputer(DFINITY)'
$ws.Range("C29").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '11.97'
$ws.Range("E29").Value = '  +6.97%  '
$ws.Range("B30").Value = 'Dai'
$ws.Range("C30").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  -0.03%  '
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.78'
$ws.Range("E31").Value = '  +10.63%  '
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.138'
$ws.Range("E32").Value = '  -4.68%  '
$ws.Range("B33").Value = 'Cronos'
$ws.Range("C33").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.184'
$ws.Range("E33").Value = '  -2.42%  '
$ws.Range("B34").Value = 'EthereumClassic'
$ws.Range("C34").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '31.08'
$ws.Range("E34").Value = '  +9.67%  '
$ws.Range("B35").Value = 'Binance-PegBSC-USD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  -0.22%  '
$ws.Range("B36").Value = 'PolygonEcosystemToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.577'
$ws.Range("E36").Value = '  +4.24%  '
$ws.Range("B37").Value = 'Bittensor'
$ws.Range("C37").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '559.53'
$ws.Range("E37").Value = '  +9.62%  '
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.81'
$ws.Range("E38").Value = '  +0.15%  '
$ws.Range("B39").Value = 'Fetch.AI'
$ws.Range("C39").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.48'
$ws.Range("E39").Value = '  -0.16%  '
$ws.Range("B40").Value = 'ARBITRUM'
$ws.Range("C40").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.945'
$ws.Range("E40").Value = '  +13.11%  '
$ws.Range("B41").Value = 'USDe'
$ws.Range("C41").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  +0.00%  '
$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.151'
$ws.Range("E42").Value = '  -0.21%  '
$ws.Range("B43").Value = 'WhiteBITCoin'
$ws.Range("C43").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '24.09'
$ws.Range("E43").Value = '  -1.15%  '
$ws.Range("B44").Value = 'ImmutableX'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.72'
$ws.Range("E44").Value = '  +2.68%  '
$ws.Range("B45").Value = 'Filecoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.73'
$ws.Range("E45").Value = '  +2.88%  '
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0418'
$ws.Range("E46").Value = '  -3.31%  '
$ws.Range("B47").Value = 'MantraDAO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.55'
$ws.Range("E47").Value = '  -2.48%  '
$ws.Range("B48").Value = 'Stacks'
$ws.Range("C48").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.19'
$ws.Range("E48").Value = '  +10.11%  '
$ws.Range("B49").Value = 'dogwifhat'
$ws.Range("C49").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.26'
$ws.Range("E49").Value = '  +4.33%  '
$ws.Range("B50").Value = 'OKB'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '53.48'
$ws.Range("E50").Value = '  +0.16%  '
$ws.Range("B51").Value = 'Cosmos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.12'
$ws.Range("E51").Value = '  -4.78%  '
